$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the two new rows with the same formatting (styles) as the rows above,
# then overwrite the values/content so the new cells keep the existing
# date / currency / integer number-format styles instead of minting new ones.
$ws.Range("A40:F41").Copy()
$ws.Range("A42:F43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 42: 四方坪站
$ws.Cells.Item(42, 1).Value = 45951
$ws.Cells.Item(42, 2).Value = "四方坪站"
$ws.Cells.Item(42, 3).Value = 9719.94
$ws.Cells.Item(42, 4).Value = 7780.35
$ws.Cells.Item(42, 5).Value = 3449.26
$ws.Cells.Item(42, 6).Value = 408

# Row 43: 高岭站
$ws.Cells.Item(43, 1).Value = 45951
$ws.Cells.Item(43, 2).Value = "高岭站"
$ws.Cells.Item(43, 3).Value = 5947.02
$ws.Cells.Item(43, 4).Value = 4690.24
$ws.Cells.Item(43, 5).Value = 1572.22
$ws.Cells.Item(43, 6).Value = 209

# Match the recorded sheet-view selection after the edit
$ws.Range("I38").Select()
